$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G403").Value = 13540.53763556403
$ws.Range("G404").Value = 27156.42071497781
$ws.Range("G405").Value = 40846.19177513176
$ws.Range("G406").Value = 54608.43209950915
$ws.Range("G407").Value = 68441.76682427636
$ws.Range("G408").Value = 82344.87003795014
$ws.Range("G409").Value = 96316.46985922658
$ws.Range("G410").Value = 110355.3534768262
$ws.Range("G411").Value = 124460.3721345508
$ws.Range("G412").Value = 138630.4460439184
$ws.Range("G413").Value = 138630.3531614433
$ws.Range("H413").Value = 13558.37666887723
$ws.Range("G414").Value = 138630.2378974515
$ws.Range("H414").Value = 27176.89939583806
$ws.Range("G415").Value = 138630.0951839034
$ws.Range("H415").Value = 40854.76868548539
$ws.Range("G416").Value = 138629.9188855589
$ws.Range("H416").Value = 54591.26941548465
$ws.Range("G417").Value = 138629.7015941503
$ws.Range("H417").Value = 68385.77508562751
$ws.Range("G418").Value = 138629.4343870651
$ws.Range("H418").Value = 82237.7519121411
$ws.Range("G419").Value = 138629.1065453013
$ws.Range("H419").Value = 96146.76274844921
$ws.Range("G420").Value = 138628.705224866
$ws.Range("H420").Value = 110112.4708128315
$ws.Range("G421").Value = 138628.2150751527
$ws.Range("H421").Value = 124134.64320263
$ws.Range("G422").Value = 138627.6177971762
$ws.Range("H422").Value = 138213.1541737743
$ws.Range("G423").Value = 138626.8916338625
$ws.Range("H423").Value = 152347.9881635155
$ws.Range("G424").Value = 138626.0107838936
$ws.Range("H424").Value = 166539.2425332463
$ws.Range("G425").Value = 138624.9447299135
$ws.Range("H425").Value = 180787.1300073083
$ws.Range("G426").Value = 138623.6574712113
$ws.Range("H426").Value = 195091.9807825598
$ws.Range("G427").Value = 138622.1066503464
$ws.Range("H427").Value = 209454.2442823869
$ws.Range("G428").Value = 138620.2425625633
$ws.Range("H428").Value = 223874.4905276532
$ws.Range("G429").Value = 138618.0070363096
$ws.Range("H429").Value = 238353.4110958576
$ws.Range("G430").Value = 138615.3321727209
$ws.Range("H430").Value = 252891.8196385478
$ws.Range("G431").Value = 138612.1389316175
$ws.Range("H431").Value = 267490.6519257569
$ws.Range("G432").Value = 138608.3355513925
$ws.Range("H432").Value = 282150.9653849542
$ws.Range("G433").Value = 138603.8157901925
$ws.Range("H433").Value = 296873.9381007565
$ws.Range("G434").Value = 138598.4569760475
$ws.Range("H434").Value = 311660.8672403846
$ws.Range("G435").Value = 138592.1178541208
$ws.Range("H435").Value = 326513.1668686886
$ws.Range("G436").Value = 138584.6362200768
$ws.Range("H436").Value = 341432.3651154226
$ws.Range("G437").Value = 138575.8263297384
$ws.Range("H437").Value = 356420.1006564795
$ws.Range("G438").Value = 138565.4760767682
$ws.Range("H438").Value = 371478.1184698872
$ws.Range("G439").Value = 138553.3439321042
$ws.Range("H439").Value = 386608.264826725
$ws.Range("G440").Value = 138539.1556413432
$ws.Range("H440").Value = 401812.4814765948
$ws.Range("G441").Value = 138522.6006792327
$ws.Range("H441").Value = 417092.7989871294
$ws.Range("G442").Value = 138503.3284639317
$ws.Range("H442").Value = 432451.3291971076
$ws.Range("G443").Value = 138480.9443377552
$ws.Range("H443").Value = 447890.256743238
$ws.Range("G444").Value = 138455.0053257407
$ws.Range("H444").Value = 463411.8296216264
$ws.Range("G445").Value = 138425.0156885694
$ws.Range("H445").Value = 479018.3487463495
$ws.Range("G446").Value = 138390.4222921351
$ws.Range("H446").Value = 494712.1564695811
$ws.Range("G447").Value = 138350.6098223472
$ws.Range("H447").Value = 510495.6240303415
$ws.Range("G448").Value = 138304.8958805603
$ws.Range("H448").Value = 526371.1379023121
$ws.Range("G449").Value = 138252.5260022586
$ws.Range("H449").Value = 542341.0850153227
$ws.Range("G450").Value = 138192.6686492432
$ws.Range("H450").Value = 558407.8368300665
$ws.Range("G451").Value = 138124.4102334549
$ws.Range("H451").Value = 574573.7322516136
$ws.Range("G452").Value = 138046.750238613
$ws.Range("H452").Value = 590841.0593741719
